$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 48.394923
$ws.Range("H2").Value = 145.184769
$ws.Range("I2").Value = 0.9646625188736886
$ws.Range("J2").Value = 0.9646625188736887
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 7.321929333333333
$ws.Range("N2").Value = 21.965788
$ws.Range("Q2").Value = 354.344206298108
$ws.Range("R2").Value = 3189.097856682972
$ws.Range("S2").Value = 0.9646625188736886
$ws.Range("T2").Value = 0.9646625188736887

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.9288356666666666
$ws.Range("H3").Value = 2.786507
$ws.Range("I3").Value = 0.01851460645626791
$ws.Range("J3").Value = 0.01851460645626791
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 7.321929333333333
$ws.Range("N3").Value = 21.965788
$ws.Range("Q3").Value = 6.800869113612888
$ws.Range("R3").Value = 61.207822022516
$ws.Range("S3").Value = 0.01851460645626791
$ws.Range("T3").Value = 0.01851460645626791

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.8439653333333333
$ws.Range("H4").Value = 2.531896
$ws.Range("I4").Value = 0.0168228746700435
$ws.Range("J4").Value = 0.0168228746700435
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 7.321929333333333
$ws.Range("N4").Value = 21.965788
$ws.Range("Q4").Value = 6.179454530449778
$ws.Range("R4").Value = 55.615090774048
$ws.Range("S4").Value = 0.0168228746700435
$ws.Range("T4").Value = 0.0168228746700435
